$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(556, 4.5, 1.43, 112.59),
    @(420, 4.6, 1.5, 88.87),
    @(420, 4.6, 1.5, 19.32),
    @(420, 4.6, 1.5, 28.98)
)

$startRow = 3
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
    $ws.Cells.Item($row, 4).Value = $data[$i][3]
}
